$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; existing rows 3-30 shift down to 4-31.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly price record.
# (Columns that repeat the same categorical values as the row below it.)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45163
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101007
$ws.Range("J3").Value = "Kiwi"
$ws.Range("K3").Value = "Hayward"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 270
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("Q3").Value = "$/bandeja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1083
$ws.Range("T3").Value = 18
